$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.453.31'
$ws.Range("E2").Value = '  +5.71%  '

$ws.Range("D3").Value = '3.389.12'
$ws.Range("E3").Value = '  +6.14%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.35'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +7.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.30'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '3.391.95'
$ws.Range("E8").Value = '  +5.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("E10").Value = '  +1.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.120'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.436'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.07%  '

$ws.Range("D13").Value = '3.974.55'
$ws.Range("E13").Value = '  +6.27%  '

$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000184'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.02'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.64%  '

$ws.Range("D17").Value = '63.529.31'
$ws.Range("E17").Value = '  +5.83%  '

$ws.Range("D18").Value = '3.397.80'
$ws.Range("E18").Value = '  +6.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.71%  '

$ws.Range("E20").Value = '  +5.20%  '

$ws.Range("E21").Value = '  +2.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.10'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.535'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.83'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.53'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +11.09%  '

$ws.Range("E27").Value = '  +18.65%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +6.00%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.50'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +6.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.12'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.80%  '

$ws.Range("E33").Value = '  +10.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.59'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.61%  '

$ws.Range("E35").Value = '  +2.38%  '

$ws.Range("E36").Value = '  +9.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.53'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.21%  '

$ws.Range("E38").Value = '  +12.04%  '

$ws.Range("E39").Value = '  +4.24%  '

$ws.Range("E40").Value = '  +6.30%  '

$ws.Range("D41").Value = '2.878.06'
$ws.Range("E41").Value = '  +1.95%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0322'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.762'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.07'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.81%  '

$ws.Range("E45").Value = '  +1.18%  '

$ws.Range("E46").Value = '  +7.96%  '

$ws.Range("D47").Value = '3.436.02'
$ws.Range("E47").Value = '  +6.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.09'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '299.40'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +12.96%  '

$ws.Range("E50").Value = '  -0.53%  '

$ws.Range("E51").Value = '  +2.84%  '
